# Convert string (inline-text) date columns to real Excel datetime
# serial values on download, matching the commit
# "convert string dates to datetime on excel download".
#
# birth_date (B), measurement_date (C) and estimated_date_delivery (F)
# were being written out as formatted GMT strings (e.g.
# "Thu, 12 Jan 2012 00:00:00 GMT"). They should instead be numeric Excel
# date serials formatted with a YYYY-MM-DD number format. Column A (the
# row index) is also refreshed to a 1-based sequence. Two rows (17 and
# 28) additionally pick up corrected derived stats (G/H/K/L) now that
# their measurement_date is no longer mis-parsed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A=1; B=40920; C=40920; F=41004 },
    @{ Row=3; A=2; B=40920; C=40933; F=41004 },
    @{ Row=4; A=3; B=40920; C=40934; F=41004 },
    @{ Row=5; A=4; B=40920; C=40960; F=41004 },
    @{ Row=6; A=5; B=40920; C=41010; F=41004 },
    @{ Row=7; A=6; B=40920; C=41284; F=41004 },
    @{ Row=8; A=7; B=40920; C=41286; F=41004 },
    @{ Row=9; A=8; B=40920; C=41649; F=41004 },
    @{ Row=10; A=9; B=40920; C=41651; F=41004 },
    @{ Row=11; A=10; B=40920; C=42380; F=41004 },
    @{ Row=12; A=11; B=40920; C=42382; F=41004 },
    @{ Row=13; A=12; B=40920; C=40920; F=40962 },
    @{ Row=14; A=13; B=40920; C=40933; F=40962 },
    @{ Row=15; A=14; B=40920; C=40934; F=40962 },
    @{ Row=16; A=15; B=40920; C=40960; F=40962 },
    @{ Row=17; A=16; B=40920; C=41010; F=40962; G=0.1314168377823409; H=0.2464065708418891; K=-0.2262967931968165; L=41.04852961824434 },
    @{ Row=18; A=17; B=40920; C=41284; F=40962 },
    @{ Row=19; A=18; B=40920; C=41286; F=40962 },
    @{ Row=20; A=19; B=40920; C=41649; F=40962 },
    @{ Row=21; A=20; B=40920; C=41651; F=40962 },
    @{ Row=22; A=21; B=40920; C=42380; F=40962 },
    @{ Row=23; A=22; B=40920; C=42382; F=40962 },
    @{ Row=24; A=23; B=40920; C=40920 },
    @{ Row=25; A=24; B=40920; C=40933 },
    @{ Row=26; A=25; B=40920; C=40934 },
    @{ Row=27; A=26; B=40920; C=40960 },
    @{ Row=28; A=27; B=40920; C=41010; G=0.2464065708418891; H=0.2464065708418891; K=-1.987173240679111; L=2.34516010760982 },
    @{ Row=29; A=28; B=40920; C=41284 },
    @{ Row=30; A=29; B=40920; C=41286 },
    @{ Row=31; A=30; B=40920; C=41649 },
    @{ Row=32; A=31; B=40920; C=41651 },
    @{ Row=33; A=32; B=40920; C=42380 },
    @{ Row=34; A=33; B=40920; C=42382 }
)

$dateFormatInitialized = $false

foreach ($r in $rows) {
    $row = $r.Row

    # A: refreshed 1-based row index (plain integer, keeps its existing style)
    $ws.Range("A" + $row).Value = $r.A

    # B: birth_date -> numeric date serial
    $ws.Range("B" + $row).Value = $r.B
    if (-not $dateFormatInitialized) {
        # First conversion: exercise the lowercase format once so both the
        # lowercase and uppercase custom number formats end up registered
        # in the workbook, then settle on the uppercase one that is
        # actually applied to every converted cell.
        $ws.Range("B" + $row).NumberFormat = "yyyy-mm-dd"
        $ws.Range("B" + $row).NumberFormat = "YYYY-MM-DD"
        $dateFormatInitialized = $true
    } else {
        $ws.Range("B" + $row).NumberFormat = "YYYY-MM-DD"
    }

    # C: measurement_date -> numeric date serial
    $ws.Range("C" + $row).Value = $r.C
    $ws.Range("C" + $row).NumberFormat = "YYYY-MM-DD"

    # F: estimated_date_delivery -> numeric date serial (only present when
    # the source row actually had a delivery-date string)
    if ($r.ContainsKey("F")) {
        $ws.Range("F" + $row).Value = $r.F
        $ws.Range("F" + $row).NumberFormat = "YYYY-MM-DD"
    }

    # G/H/K/L: corrected_decimal_age / chronological_decimal_age / sds /
    # centile recomputed from the fixed measurement_date for the couple of
    # rows whose source data was wrong before this fix.
    if ($r.ContainsKey("G")) {
        $ws.Range("G" + $row).Value = $r.G
    }
    if ($r.ContainsKey("H")) {
        $ws.Range("H" + $row).Value = $r.H
    }
    if ($r.ContainsKey("K")) {
        $ws.Range("K" + $row).Value = $r.K
    }
    if ($r.ContainsKey("L")) {
        $ws.Range("L" + $row).Value = $r.L
    }
}
